# "Ajout fonction score sur question"
# Adds a new "id" column (1..170) to the end of the Tableau1 table on the
# single worksheet, numbering each question row, and leaves the view
# positioned on the last-edited cell (K169), matching the author's
# final cursor position after extending the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with a new trailing column -------------------------
# ListColumns.Add() appends a column right after the table's last column
# and grows both the table ref and the worksheet's autofilter/dimension to
# match (A1:I171 -> A1:J171).
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add()

# Name the new header cell "id" (this is what actually renames the
# corresponding table column).
$ws.Range("J1").Value = "id"

# --- Fill the new column with a running question index (1..170) ---------
for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 10).Value = $r - 1
}

# --- Restore the cursor/selection to where the author left it -----------
[void]$excel.Goto($ws.Range("B133"), $true)
[void]$ws.Range("K169").Select()
